$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.756210803985596
$ws.Range("B1").Value = 3.133386373519897
$ws.Range("C1").Value = 2.863378286361694
$ws.Range("D1").Value = 3.244951248168945
$ws.Range("E1").Value = 2.231599569320679
